$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$changes = @{
    "G2"  = 2.25
    "H2"  = 3.15
    "I2"  = 3.15
    "L2"  = 3.7
    "N2"  = 7.6
    "S2"  = 1.38
    "T2"  = 2.82
    "U2"  = 1.57
    "V2"  = 2.25
    "W2"  = 9.5
    "X2"  = 13
    "Z2"  = 24
    "AB2" = 21
    "AC2" = 7.6
    "AD2" = 6.3
    "AG2" = 300
    "AH2" = 10.5
    "AI2" = 17.5
    "AJ2" = 10.75
    "AK2" = 45
    "AM2" = 29
    "AN2" = 4.3
    "AP2" = 17
    "AR2" = 65
    "AS2" = 175
    "AT2" = 2.82
    "AW2" = 5.2
    "AX2" = 18
    "AY2" = 23
    "AZ2" = 90
    "BA2" = 110
    "BB2" = 250
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
